# "Creacion de pdf con imagenes"
#
# The LOGIN worksheet is rebuilt: its lead blank/header row is removed so the
# USERNAME/PASSWORD header moves to row 1 and the credential values move to
# row 2 (dropping the stray styled D/E cells that used to sit next to them).
# Rebuilding the sheet (delete + re-add) is what naturally reproduces the new
# internal sheetId (8) and drops the old codeName the diff shows disappearing.
# The MAIN sheet only has its remembered cursor/selection cell changed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# MAIN: just a cursor/selection move (E20 -> B33), no data changes.
# ---------------------------------------------------------------------
$main = $wb.Worksheets.Item("MAIN")
$main.Range("B33").Select()

# ---------------------------------------------------------------------
# LOGIN: drop the old sheet and rebuild it in the same slot (between MAIN
# and ITEM) so the credentials sit on rows 1-2 instead of 2-3.
# ---------------------------------------------------------------------
$oldLogin = $wb.Worksheets.Item("LOGIN")
$oldLogin.Delete()

$itemSheet = $wb.Worksheets.Item("ITEM")
$login = $wb.Worksheets.Add($itemSheet)
$login.Name = "LOGIN"

# Header row (was row 2, now row 1)
$login.Range("A1").Value = "USERNAME"
$login.Range("B1").Value = "PASSWORD"

# Credential row (was row 3, now row 2) - keep the hyperlink styling on A2
$login.Range("A2").Value = "black_panther_xaa@hotmail.com"
$login.Range("B2").Value = "futbol01"

# Restore the mailto hyperlink, now anchored on A2 instead of A3, then
# (re)apply the "Hyperlink" cell style so A2 keeps the underlined/colored look
$login.Hyperlinks.Add($login.Range("A2"), "mailto:black_panther_xaa@hotmail.com")
$login.Range("A2").Style = "Hyperlink"

# Match the original column widths as closely as this model allows
$login.Columns.Item(1).ColumnWidth = 30.5
$login.Columns.Item(2).ColumnWidth = 10.333333333333334

# Final cursor/selection on the rebuilt sheet
$login.Range("A30").Select()
